$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("plotConfiguration")
$ws.Range("I1").Value = "xValuesLimits"
$ws.Range("J1").Value = "yValuesLimits"
$ws.Range("J2").Select()
